$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3, shifting the existing data rows (previously 3..85) down to 4..86.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new data record.
$ws.Range("A3").Value2 = 1
$ws.Range("B3").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value2 = "Arica y Parinacota"
$ws.Range("D3").Value2 = 44691
$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("E3").Value2 = 15
$ws.Range("F3").Value2 = 100112038
$ws.Range("G3").Value2 = "Cebollín baby"
$ws.Range("H3").Value2 = "Sin especificar"
$ws.Range("I3").Value2 = "Primera"
$ws.Range("J3").Value2 = 270
$ws.Range("K3").Value2 = 2000
$ws.Range("L3").Value2 = 2500
$ws.Range("M3").Value2 = 2250
$ws.Range("N3").Value2 = "`$/paquete 1,5 a 2 kilos"
$ws.Range("O3").Value2 = "Región de Arica y Parinacota"
$ws.Range("P3").Value2 = 1125
$ws.Range("Q3").Value2 = 2
$ws.Range("R3").Value2 = "Hortaliza"
